# Apply the May 9th gyroscope data changes:
# 1. Insert 7 new rows of data right after the header row (new rows 2-8),
#    pushing the existing data rows down.
# 2. Append 3 new rows of data at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 7 rows after row 1 (the header) ---
$insertRange = $ws.Range("A2:C8")
$insertRange.EntireRow.Insert()
$insertRange.EntireRow.ClearFormats()

$newTopRows = @(
    @(-0.2191115766763687, -0.8925584554672241, 0.058589544147253),
    @(-0.4737097918987274, -0.3764378130435943, -0.0346210934221744),
    @(2.269346237182617, 2.84225869178772, -1.580319762229919),
    @(2.80038046836853, -2.297109603881836, -0.6647250056266785),
    @(-2.487059593200684, -0.9021458625793456, 2.389388084411621),
    @(-9.80595874786377, 2.084323167800904, 0.3504720032215118),
    @(1.612078070640564, 6.0146164894104, 2.238120555877685)
)

for ($i = 0; $i -lt $newTopRows.Count; $i++) {
    $rowNum = 2 + $i
    $rowData = $newTopRows[$i]
    $ws.Cells.Item($rowNum, 1).Value = $rowData[0]
    $ws.Cells.Item($rowNum, 2).Value = $rowData[1]
    $ws.Cells.Item($rowNum, 3).Value = $rowData[2]
}

# --- Step 2: append 3 new rows at the end (rows 29-31) ---
$newBottomRows = @(
    @(-7.232280254364014, -3.663843870162964, 11.63588333129883),
    @(1.14389431476593, 9.762216567993164, -6.083192825317383),
    @(2.391319036483765, -6.093712329864502, 2.177400588989258)
)

for ($i = 0; $i -lt $newBottomRows.Count; $i++) {
    $rowNum = 29 + $i
    $rowData = $newBottomRows[$i]
    $ws.Cells.Item($rowNum, 1).Value = $rowData[0]
    $ws.Cells.Item($rowNum, 2).Value = $rowData[1]
    $ws.Cells.Item($rowNum, 3).Value = $rowData[2]
}
